# Update functions and Data Model (#50)
# Adds a new "Authorship Resource" column (I) to Sheet1, populated with the
# resource's authorship attribution, and updates the selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I
$ws.Range("I1").Value = "Authorship Resource"

# New values for column I, rows 2-5
$ws.Range("I2").Value = "Daniela Subotic, Noémi Villars-Amberg"
$ws.Range("I3").Value = "Daniela Subotic, Noémi Villars-Amberg"
$ws.Range("I4").Value = "Daniela Subotic, Noémi Villars-Amberg"
$ws.Range("I5").Value = "Daniela Subotic, Noémi Villars-Amberg"

# Reflect the new active selection in the sheet view
$ws.Range("J15").Select()
